$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference cell (column C) used to copy a guaranteed "General" number format
# back onto column D cells after forcing them to Text, so price strings like
# "26.702.87" / "1.531.32" (which are not valid numbers anyway) and plain
# decimals like "205.54" are both stored as literal text, matching the source data.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.702.87"
$ws.Range("D2").NumberFormat = $ws.Range("C2").NumberFormat
$ws.Range("E2").Value = "  -0.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.531.32"
$ws.Range("D3").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("E3").Value = "  -1.94%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.54"
$ws.Range("D5").NumberFormat = $ws.Range("C5").NumberFormat
$ws.Range("E5").Value = "  -0.36%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("E9").Value = "  -1.03%  "
$ws.Range("E10").Value = "  -0.63%  "
$ws.Range("E11").Value = "  -1.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.749.25"
$ws.Range("D12").NumberFormat = $ws.Range("C12").NumberFormat
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.530.20"
$ws.Range("D13").NumberFormat = $ws.Range("C13").NumberFormat
$ws.Range("E13").Value = "  -2.15%  "
$ws.Range("E14").Value = "  -1.50%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.506"
$ws.Range("D15").NumberFormat = $ws.Range("C15").NumberFormat
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "26.686.62"
$ws.Range("D16").NumberFormat = $ws.Range("C16").NumberFormat
$ws.Range("E16").Value = "  -0.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.05"
$ws.Range("D17").NumberFormat = $ws.Range("C17").NumberFormat
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "212.18"
$ws.Range("D18").NumberFormat = $ws.Range("C18").NumberFormat
$ws.Range("E19").Value = "  +0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.19"
$ws.Range("D20").NumberFormat = $ws.Range("C20").NumberFormat
$ws.Range("E20").Value = "  -2.21%  "
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.00"
$ws.Range("D22").NumberFormat = $ws.Range("C22").NumberFormat
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("E23").Value = "  -3.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.95"
$ws.Range("D24").NumberFormat = $ws.Range("C24").NumberFormat
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.86"
$ws.Range("D25").NumberFormat = $ws.Range("C25").NumberFormat
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.58"
$ws.Range("D26").NumberFormat = $ws.Range("C26").NumberFormat
$ws.Range("E26").Value = "  -2.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.79"
$ws.Range("D27").NumberFormat = $ws.Range("C27").NumberFormat
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E28").Value = "  -0.05%  "
$ws.Range("E29").Value = "  -1.36%  "
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("E31").Value = "  -2.01%  "
$ws.Range("E32").Value = "  +2.78%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.360.13"
$ws.Range("D33").NumberFormat = $ws.Range("C33").NumberFormat
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -0.39%  "
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.941"
$ws.Range("D36").NumberFormat = $ws.Range("C36").NumberFormat
$ws.Range("E36").Value = "  +1.07%  "
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.521"
$ws.Range("D39").NumberFormat = $ws.Range("C39").NumberFormat
$ws.Range("E39").Value = "  +0.29%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.71"
$ws.Range("D40").NumberFormat = $ws.Range("C40").NumberFormat
$ws.Range("E40").Value = "  +6.40%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").NumberFormat = $ws.Range("C41").NumberFormat
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("E42").Value = "  -1.75%  "
$ws.Range("E43").Value = "  +0.11%  "
$ws.Range("E44").Value = "  +1.03%  "
$ws.Range("E45").Value = "  -1.45%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "62.43"
$ws.Range("D46").NumberFormat = $ws.Range("C46").NumberFormat
$ws.Range("E46").Value = "  -1.21%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.663.44"
$ws.Range("D47").NumberFormat = $ws.Range("C47").NumberFormat
$ws.Range("E47").Value = "  -1.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.30"
$ws.Range("D48").NumberFormat = $ws.Range("C48").NumberFormat
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("D49").NumberFormat = $ws.Range("C49").NumberFormat
$ws.Range("E49").Value = "  +1.93%  "
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0942"
$ws.Range("D51").NumberFormat = $ws.Range("C51").NumberFormat
$ws.Range("E51").Value = "  -0.59%  "
